$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Relative frequency [%] (col B) and S.D. (col C) values for rows 12-72
$bValues = @(
  0, 0, 0.095590680837631226, 0.60643035173416138, 1.4987994432449341, 2.5590324401855469,
  3.5279166698455811, 4.19200325012207, 4.4326229095458984, 4.2379312515258789, 3.6868877410888672,
  2.9160935878753662, 2.0806283950805664, 1.3184617757797241, 0.725055456161499, 0.341017484664917,
  0.15187500417232513, 0.0958925187587738, 0.10026177763938904, 0.14344993233680725,
  0.19715073704719543, 0.23501701653003693, 0.2598978579044342, 0.26955568790435791,
  0.26068574190139771, 0.23236016929149628, 0.18875923752784729, 0.14033815264701843,
  0.10195136815309525, 0.08616231381893158, 0.090176895260810852, 0.097623430192470551,
  0.33288782835006714, 0.78972560167312622, 1.3259538412094116, 1.755842924118042, 1.9286437034606934,
  1.795910120010376, 1.4694110155105591, 1.1720955371856689, 1.061733603477478, 1.1856088638305664,
  1.4081041812896729, 1.4190545082092285, 1.212550163269043, 0.86168420314788818, 0.49653580784797668,
  0.25098052620887756, 0.18041279911994934, 0.15623000264167786, 0.18230223655700684,
  0.38992920517921448, 0.88348358869552612, 1.5032562017440796, 1.9889378547668457,
  2.1215076446533203, 1.8239263296127319, 1.2030737400054932, 0.51888000965118408, 0.200624480843544,
  0.53621417284011841
)

$cValues = @(
  0, 0, 0.071805007755756378, 0.097744680941104889, 0.40017604827880859, 0.84054863452911377,
  1.2723718881607056, 1.5764970779418945, 1.6798771619796753, 1.5659104585647583, 1.2697107791900635,
  0.86359214782714844, 0.44081589579582214, 0.14586842060089111, 0.24480704963207245,
  0.32096239924430847, 0.2630552351474762, 0.1247498095035553, 0.11439704149961472,
  0.24846257269382477, 0.34147509932518005, 0.37291008234024048, 0.33979034423828125,
  0.2768130898475647, 0.22927513718605042, 0.2033819854259491, 0.16767686605453491,
  0.12195190787315369, 0.11904421448707581, 0.14923751354217529, 0.15619096159934998,
  0.13324242830276489, 0.42184650897979736, 1.2731107473373413, 2.2653806209564209,
  3.0412089824676514, 3.3185341358184814, 2.9792938232421875, 2.0798029899597168, 0.843402624130249,
  0.36177638173103333, 1.10057532787323, 1.3246496915817261, 1.3153784275054932, 1.0802453756332397,
  0.74770873785018921, 0.52241355180740356, 0.43214708566665649, 0.31248411536216736,
  0.27059829235076904, 0.31575673818588257, 0.37382104992866516, 0.83384084701538086,
  1.6747992038726807, 2.4209990501403809, 2.7135639190673828, 2.4086086750030518, 1.6297746896743774,
  0.72560840845108032, 0.18811064958572388, 0.92875027656555176
)

$startRow = 12
for ($i = 0; $i -lt $bValues.Count; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 2).Value = $bValues[$i]
  $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

Write-Output "Updated B12:C72 with new Relative frequency and S.D. values"